$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the timestamp on the existing last row (A26): tiny precision fix
# from the same day's data re-retrieval.
$ws.Cells.Item(26, 1).Value = 44339.77832643634

# Append the newly retrieved row (27) of job numbers.
$ws.Cells.Item(27, 1).Value = 44340.78333358185
$ws.Cells.Item(27, 2).Value = 73767
$ws.Cells.Item(27, 3).Value = 62102
$ws.Cells.Item(27, 4).Value = 3313
$ws.Cells.Item(27, 5).Value = 2084
$ws.Cells.Item(27, 6).Value = 1472
$ws.Cells.Item(27, 7).Value = 19204
$ws.Cells.Item(27, 8).Value = 1360
$ws.Cells.Item(27, 9).Value = 831
$ws.Cells.Item(27, 10).Value = 205

# Carry the date/time number format from A26 onto the new A27 cell so the
# timestamp renders the same way (style index 2 in the sheet).
$ws.Cells.Item(27, 1).NumberFormat = $ws.Cells.Item(26, 1).NumberFormat
